# ScenTrade__Trade_Links.xlsx edit
#
# The underlying data change in this commit is a rename of the regions
# used on the "BI" sheet from the Arizona naming (AZ1/AZ2/AZ3) to the
# German naming (DE1/DE2/DE3) - this workbook was repurposed from the
# TIMES_AZ model to the TIMES-DE model. Everything else in the original
# diff (fileVersion/build numbers, the author's local folder path,
# window geometry, customXml part re-numbering, x14ac:dyDescent / row
# height tweaks, ...) is incidental save noise produced by opening and
# re-saving the file on a different machine/Excel build - the author's
# own commit message says as much ("Changes hvor jeg ikke har gjort
# noget" = "Changes where I haven't done anything"). That noise isn't
# something a user performs through the Excel UI/object model, so this
# script focuses on reproducing the real edit: the AZ1/AZ2/AZ3 -> DE1/
# DE2/DE3 rename, plus leaving the workbook with the "BI" sheet active
# and cell I11 selected, which is the other observable, intentional
# state captured in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BI")

# Header row (row 4) labels
$ws.Range("D4").Value = "DE1"
$ws.Range("E4").Value = "DE2"
$ws.Range("F4").Value = "DE3"

# Row labels (column C, rows 5-7) mirror the header row
$ws.Range("C5").Value = "DE1"
$ws.Range("C6").Value = "DE2"
$ws.Range("C7").Value = "DE3"

# Leave "BI" as the active sheet/tab with I11 selected, matching the
# saved view state in the updated workbook.
$ws.Activate()
$ws.Range("I11").Select()
